$wb = $excel.ActiveWorkbook

# 1. Rename shared string "elc_win-CHE" -> "elc_won-CHE" across the workbook
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("elc_win-CHE", "elc_won-CHE")
}

# 2. Update the two long timeslice-list strings on ev_charging_uc (C13, C14)
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("C13").Value = "S1aH3,S2e0416h11,S2e0416h18,S2i0429h11,S2i0429h12,S2i0429h15,S2j0518h16,S2j0518h18,S2k0531h17,S4l0918h10,S4l0918h12,S4l0918h13,S4l0918h17,S4m1019h11,S4n1020h15,S5p1104h15,S1aH4,S1b0201h08,S2aH2,S2f0420h17,S2i0429h16,S4aH6,S4l0918h18,S4n1020h12,S4n1020h18,S4o1023h15,S1aH6,S1b0201h11,S1b0201h14,S1b0201h16,S1b0201h17,S1b0201h18,S2g0421h11,S2g0421h12,S2g0421h14,S2h0427h12,S2h0427h15,S2j0518h13,S2k0531h13,S2k0531h16,S4l0918h09,S4l0918h11,S5p1104h10,S5p1104h18,S1b0201h13,S1b0201h15,S2e0416h10,S2e0416h14,S2e0416h17,S2f0420h10,S2i0429h17,S2k0531h14,S4m1019h12,S4n1020h08,S4n1020h11,S1b0201h10,S1b0201h12,S2d0415h08,S2d0415h10,S2d0415h13,S2e0416h09,S2f0420h14,S2i0429h18,S2j0518h14,S4aH2,S4m1019h14,S4o1023h07,S4o1023h16,S1c0205h15,S2f0420h07,S2g0421h09,S2h0427h10,S2j0518h10,S2k0531h11,S3aH3,S4n1020h07,S5p1104h14,S1aH5,S1c0205h16,S2aH4,S2aH5,S2d0415h11,S2d0415h17,S2e0416h16,S2g0421h08,S2g0421h13,S2h0427h07,S2h0427h11,S2k0531h10,S4m1019h09,S4m1019h18,S4o1023h17,S5aH2,S5aH5,S5p1104h16,S6aH2,S1b0201h09,S1c0205h17,S2aH6,S2e0416h15,S2f0420h13,S2g0421h10,S2g0421h15,S2h0427h18,S2j0518h09,S2j0518h12,S2k0531h08,S3aH5,S4m1019h13,S4m1019h15,S4o1023h09,S5aH3,S5p1104h08,S1aH2,S1aH7,S1c0205h10,S2aH3,S2e0416h07,S2f0420h09,S2f0420h18,S2i0429h08,S2j0518h11,S2j0518h15,S2j0518h17,S2k0531h18,S3aH6,S4aH7,S4l0918h08,S4o1023h11,S5aH6,S5aH7,S6aH6,S1b0201h07,S2d0415h12,S2g0421h16,S2j0518h08,S4aH5,S4l0918h16,S5p1104h07,S1c0205h09,S1c0205h11,S2aH7,S2d0415h09,S2h0427h16,S2k0531h15,S3aH4,S4aH3,S4aH4,S4m1019h10,S4m1019h16,S5p1104h09,S5p1104h11,S1c0205h13,S2d0415h07,S2f0420h08,S2f0420h11,S2h0427h09,S2h0427h14,S2i0429h13,S3aH7,S4l0918h07,S4l0918h14,S4o1023h18,S2d0415h16,S2e0416h13,S2f0420h12,S2f0420h15,S2g0421h17,S2h0427h08,S2j0518h07,S2k0531h12,S3aH2,S4o1023h14,S5p1104h12,S1c0205h07,S1c0205h08,S2d0415h14,S2d0415h15,S2i0429h07,S2i0429h09,S2i0429h14,S4n1020h10,S4n1020h13,S4n1020h14,S4n1020h17,S4o1023h08,S4o1023h12,S6aH5,S1c0205h18,S2d0415h18,S2e0416h08,S2f0420h16,S4l0918h15,S5aH4,S5p1104h13,S5p1104h17,S6aH3,S6aH4,S1c0205h12,S1c0205h14,S2e0416h12,S2g0421h07,S2g0421h18,S2h0427h13,S2h0427h17,S2i0429h10,S2k0531h07,S2k0531h09,S4m1019h07,S4m1019h08,S4m1019h17,S4n1020h09,S4n1020h16,S4o1023h10,S4o1023h13,S6aH7"
$wsEv.Range("C14").Value = "S1b0201h02,S2d0415h03,S2d0415h20,S2g0421h05,S2g0421h23,S2h0427h20,S2h0427h21,S2i0429h20,S2j0518h01,S2j0518h19,S2j0518h23,S4l0918h06,S4m1019h20,S4n1020h19,S1b0201h20,S1b0201h24,S2d0415h05,S2f0420h01,S2h0427h01,S4aH1,S4l0918h24,S4n1020h02,S4o1023h19,S6aH1,S1b0201h19,S1c0205h06,S2d0415h21,S2g0421h04,S2i0429h04,S4l0918h20,S4m1019h02,S4n1020h21,S4o1023h05,S5p1104h22,S6aH8,S1aH1,S1b0201h23,S2d0415h22,S2f0420h21,S2g0421h24,S2i0429h03,S4n1020h03,S4n1020h20,S2aH1,S2g0421h02,S2g0421h03,S2h0427h05,S2h0427h06,S2i0429h19,S2k0531h19,S4l0918h01,S4m1019h22,S4n1020h23,S4o1023h20,S4o1023h24,S5aH8,S5p1104h03,S1b0201h01,S2e0416h01,S2f0420h20,S2g0421h21,S4l0918h22,S5p1104h19,S1c0205h04,S1c0205h19,S1c0205h22,S2f0420h04,S2j0518h05,S4l0918h23,S4m1019h19,S4n1020h04,S4n1020h06,S1aH8,S1b0201h22,S1c0205h03,S2d0415h23,S2g0421h01,S2g0421h20,S2i0429h01,S2j0518h04,S2j0518h21,S4m1019h06,S4o1023h04,S5aH1,S5p1104h05,S1b0201h03,S1b0201h06,S1c0205h02,S1c0205h05,S2d0415h24,S2e0416h20,S2e0416h23,S2f0420h03,S2h0427h03,S2h0427h23,S2i0429h21,S2j0518h06,S2k0531h05,S4l0918h21,S4m1019h23,S4n1020h24,S4o1023h02,S4o1023h21,S5p1104h01,S1c0205h23,S2e0416h02,S2e0416h04,S2f0420h22,S2h0427h22,S2i0429h24,S2j0518h03,S2j0518h24,S4n1020h01,S5p1104h02,S1b0201h05,S1c0205h01,S1c0205h20,S2f0420h05,S2f0420h24,S2g0421h06,S2h0427h19,S2j0518h20,S2k0531h06,S4m1019h04,S4o1023h03,S5p1104h20,S2aH8,S2e0416h19,S2f0420h19,S2g0421h22,S2i0429h02,S2i0429h22,S2i0429h23,S2k0531h04,S2k0531h24,S4l0918h02,S4l0918h05,S4m1019h24,S4n1020h05,S4n1020h22,S5p1104h04,S5p1104h24,S1c0205h24,S2d0415h04,S2h0427h04,S2k0531h21,S3aH8,S4m1019h01,S5p1104h21,S1b0201h04,S2e0416h03,S2e0416h05,S2e0416h06,S2e0416h22,S2f0420h02,S2i0429h05,S2k0531h20,S2k0531h22,S4aH8,S4l0918h03,S4o1023h01,S4o1023h23,S5p1104h06,S5p1104h23,S2d0415h01,S2d0415h02,S2d0415h06,S2d0415h19,S2g0421h19,S2h0427h24,S2i0429h06,S2k0531h01,S2k0531h23,S4l0918h04,S4m1019h03,S4m1019h05,S4m1019h21,S1b0201h21,S1c0205h21,S2e0416h21,S2e0416h24,S2f0420h06,S2f0420h23,S2h0427h02,S2j0518h02,S2j0518h22,S2k0531h02,S2k0531h03,S3aH1,S4l0918h19,S4o1023h06,S4o1023h22"

# 3. Update the seasonal split values on re_profiles (M/N columns, rows 4-9)
$wsRe = $wb.Worksheets.Item("re_profiles")
$wsRe.Range("N4").Value = 0.080605632899210897
$wsRe.Range("M5").Value = "S3"
$wsRe.Range("N5").Value = 0.26179508564829668
$wsRe.Range("N6").Value = 0.19493167383075641
$wsRe.Range("M7").Value = "S2"
$wsRe.Range("N7").Value = 0.40660807082825434
$wsRe.Range("M8").Value = "S1"
$wsRe.Range("N8").Value = 0.17653172515557838
$wsRe.Range("N9").Value = 0.079527811637903401
